$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3: One for the Books | Leather Grimoire
$ws.Range("H3").Value = 22400
$ws.Range("J3").Value = 22400
$ws.Range("L3").Value = 22400
$ws.Range("N3").Value = -22628

# Row 17: One for the Road | Potion
$ws.Range("H17").Value = 563038.2
$ws.Range("J17").Value = 600532.5600000001
$ws.Range("L17").Value = 1801597.68
$ws.Range("N17").Value = -1801933.68

# Row 18: You Grow, Girl | Growth Formula Beta
$ws.Range("H18").Value = 1190
$ws.Range("J18").Value = 590
$ws.Range("L18").Value = 590
$ws.Range("N18").Value = -1158

# Row 19: Unbreak My Heart | Roof Tile
$ws.Range("H19").Value = 588.8889
$ws.Range("I19").Value = 683.6667
$ws.Range("J19").Value = 541.5
$ws.Range("K19").Value = 683.6667
$ws.Range("L19").Value = 541.5
$ws.Range("M19").Value = -508.6667
$ws.Range("N19").Value = -891.5

# Row 74: Adhesive of Antipathy | Wing Glue
$ws.Range("H74").Value = 3705.348
$ws.Range("I74").Value = 3616.8125
$ws.Range("J74").Value = 3907.7144
$ws.Range("K74").Value = 3616.8125
$ws.Range("L74").Value = 3907.7144
$ws.Range("M74").Value = -2680.8125
$ws.Range("N74").Value = -5779.7144

# Row 77: It's Gonna Grow Back (L) | Wing Glue
$ws.Range("H77").Value = 3705.348
$ws.Range("I77").Value = 3616.8125
$ws.Range("J77").Value = 3907.7144
$ws.Range("K77").Value = 18084.0625
$ws.Range("L77").Value = 19538.572
$ws.Range("M77").Value = -13404.0625
$ws.Range("N77").Value = -28898.572

# Row 87: There Was a Late Fee | Noble Gold
$ws.Range("H87").Value = 33200
$ws.Range("J87").Value = 33200
$ws.Range("L87").Value = 33200
$ws.Range("N87").Value = -35696

# Row 90: A Gate Arcane Is Dragon's Bane (L) | Noble Gold
$ws.Range("H90").Value = 33200
$ws.Range("J90").Value = 33200
$ws.Range("L90").Value = 99600
$ws.Range("N90").Value = -112080

# Row 102: Spell-rebound | Marid Leather Grimoire
$ws.Range("H102").Value = 22400
$ws.Range("J102").Value = 22400
$ws.Range("L102").Value = 22400
$ws.Range("N102").Value = -28890

# Row 121: Mindful Medicine | Tincture of Mind
$ws.Range("H121").Value = 1259.1666
$ws.Range("J121").Value = 1500
$ws.Range("L121").Value = 4500
$ws.Range("N121").Value = -7994

# Row 127: Liquid Competence | Competent Craftsman's Draught
$ws.Range("H127").Value = 973.05554
$ws.Range("I127").Value = 533.125
$ws.Range("J127").Value = 1325
$ws.Range("K127").Value = 1599.375
$ws.Range("L127").Value = 3975
$ws.Range("M127").Value = 3360.625
$ws.Range("N127").Value = -13895

# Row 129: Practical Command | Commanding Craftsman's Draught
$ws.Range("H129").Value = 1144.8918
$ws.Range("I129").Value = 600
$ws.Range("J129").Value = 1160.0278
$ws.Range("K129").Value = 1800
$ws.Range("L129").Value = 3480.0834
$ws.Range("M129").Value = 3200
$ws.Range("N129").Value = -13480.0834

# Row 131: Mindful Study | Grade 5 Tincture of Mind
$ws.Range("H131").Value = 1236.1875
$ws.Range("I131").Value = 875.3077
$ws.Range("J131").Value = 2800
$ws.Range("K131").Value = 2625.9231
$ws.Range("L131").Value = 8400
$ws.Range("M131").Value = 2414.0769
$ws.Range("N131").Value = -18480

# Row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 5988.1377
$ws.Range("I141").Value = 6329
$ws.Range("J141").Value = 4352
$ws.Range("K141").Value = 18987
$ws.Range("L141").Value = 13056
$ws.Range("M141").Value = -13807
$ws.Range("N141").Value = -23416

$ws = $wb.Worksheets.Item("ARM")
# Row 80: A Squire to Inspire | Titanium Hoplon
$ws.Range("H80").Value = 32597.691
$ws.Range("I80").Value = 20000
$ws.Range("J80").Value = 33647.5
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 33647.5
$ws.Range("M80").Value = -19002
$ws.Range("N80").Value = -35643.5

# Row 83: All's Fair in Highborn Assassination (L) | Titanium Hoplon
$ws.Range("H83").Value = 32597.691
$ws.Range("I83").Value = 20000
$ws.Range("J83").Value = 33647.5
$ws.Range("K83").Value = 60000
$ws.Range("L83").Value = 100942.5
$ws.Range("M83").Value = -55008
$ws.Range("N83").Value = -110926.5

# Row 101: Art Imitates Life | Doman Steel Tabard of Fending
$ws.Range("H101").Value = 39087.43
$ws.Range("J101").Value = 39087.43
$ws.Range("L101").Value = 39087.43
$ws.Range("N101").Value = -45577.43

# Row 139: Backing up My Words | Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 86950
$ws.Range("J139").Value = 86950
$ws.Range("L139").Value = 86950
$ws.Range("N139").Value = -97230

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 6219.1304
$ws.Range("I134").Value = 7374.1177
$ws.Range("K134").Value = 22122.3531
$ws.Range("M134").Value = -19587.3531

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof | Ash Lumber
$ws.Range("H16").Value = 757.625
$ws.Range("I16").Value = 743.5
$ws.Range("K16").Value = 743.5
$ws.Range("M16").Value = -456.5

# Row 33: Tools for the Tools | Silver Battle Fork
$ws.Range("H33").Value = 7069.875
$ws.Range("I33").Value = 1312.4
$ws.Range("J33").Value = 16665.666
$ws.Range("K33").Value = 1312.4
$ws.Range("L33").Value = 16665.666
$ws.Range("M33").Value = -933.4000000000001
$ws.Range("N33").Value = -17423.666

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 1013.549
$ws.Range("I58").Value = 1001.5814
$ws.Range("J58").Value = 1077.875
$ws.Range("K58").Value = 1001.5814
$ws.Range("L58").Value = 1077.875
$ws.Range("M58").Value = -798.5814
$ws.Range("N58").Value = -1483.875

# Row 99: O Pine | Pine Lumber
$ws.Range("H99").Value = 1478.7222
$ws.Range("I99").Value = 1432.4615
$ws.Range("J99").Value = 1599
$ws.Range("K99").Value = 1432.4615
$ws.Range("L99").Value = 1599
$ws.Range("M99").Value = 65.53850000000011
$ws.Range("N99").Value = -4595

# Row 106: With a Bow on Top | Zelkova Longbow
$ws.Range("H106").Value = 32671
$ws.Range("J106").Value = 32671
$ws.Range("L106").Value = 32671
$ws.Range("N106").Value = -35195

# Row 113: Patient Patients | White Ash Lumber
$ws.Range("H113").Value = 757.625
$ws.Range("I113").Value = 743.5
$ws.Range("K113").Value = 743.5
$ws.Range("M113").Value = 1426.5

# Row 120: Kindling the Flame | Lignum Vitae Ring
$ws.Range("H120").Value = 50326
$ws.Range("J120").Value = 50326
$ws.Range("L120").Value = 50326
$ws.Range("N120").Value = -57584

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Range("H126").Value = 1478.7222
$ws.Range("I126").Value = 1432.4615
$ws.Range("J126").Value = 1599
$ws.Range("K126").Value = 4297.3845
$ws.Range("L126").Value = 4797
$ws.Range("M126").Value = -1827.3845
$ws.Range("N126").Value = -9737

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 1013.549
$ws.Range("I136").Value = 1001.5814
$ws.Range("J136").Value = 1077.875
$ws.Range("K136").Value = 3004.7442
$ws.Range("L136").Value = 3233.625
$ws.Range("M136").Value = -454.7442000000001
$ws.Range("N136").Value = -8333.625

$ws = $wb.Worksheets.Item("CUL")
# Row 55: Pagan Pastries | Pastry Fish
$ws.Range("H55").Value = 2773.6843
$ws.Range("I55").Value = 1800
$ws.Range("J55").Value = 2888.2354
$ws.Range("K55").Value = 5400
$ws.Range("L55").Value = 8664.706200000001
$ws.Range("M55").Value = -5223
$ws.Range("N55").Value = -9018.706200000001

# Row 64: The Aroma of Faith | Baked Onion Soup
$ws.Range("H64").Value = 3365.6
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 4942.6665
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 14827.9995
$ws.Range("M64").Value = -2730
$ws.Range("N64").Value = -15367.9995

# Row 67: Soup's On (L) | Baked Onion Soup
$ws.Range("H67").Value = 3365.6
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 4942.6665
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 14827.9995
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = -16699.9995

# Row 118: Teetotally | Masala Chai
$ws.Range("H118").Value = 2269.4688
$ws.Range("I118").Value = 732.5
$ws.Range("J118").Value = 2624.1538
$ws.Range("K118").Value = 2197.5
$ws.Range("L118").Value = 7872.4614
$ws.Range("M118").Value = -954.5
$ws.Range("N118").Value = -10358.4614

# Row 121: A Cookie for Your Troubles | Coffee Biscuit
$ws.Range("H121").Value = 1480.2354
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 1690.2858
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 5070.857400000001
$ws.Range("M121").Value = -190
$ws.Range("N121").Value = -7690.857400000001

# Row 129: Comfort Food | Yakow Moussaka
$ws.Range("H129").Value = 13889601
$ws.Range("I129").Value = 526
$ws.Range("J129").Value = 37038060
$ws.Range("K129").Value = 1578
$ws.Range("L129").Value = 111114180
$ws.Range("M129").Value = 3422
$ws.Range("N129").Value = -111124180

# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 1575.25
$ws.Range("J132").Value = 1916.5
$ws.Range("L132").Value = 17248.5
$ws.Range("N132").Value = -22308.5

# Row 134: Don't Knock It Till You've Tried It | Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 4203.9
$ws.Range("I134").Value = 907.8
$ws.Range("K134").Value = 2723.4
$ws.Range("M134").Value = 2346.6

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 3044.0908
$ws.Range("I80").Value = 2151.25
$ws.Range("J80").Value = 3554.2856
$ws.Range("K80").Value = 2151.25
$ws.Range("L80").Value = 3554.2856
$ws.Range("M80").Value = -1153.25
$ws.Range("N80").Value = -5550.2856

# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 3044.0908
$ws.Range("I83").Value = 2151.25
$ws.Range("J83").Value = 3554.2856
$ws.Range("K83").Value = 10756.25
$ws.Range("L83").Value = 17771.428
$ws.Range("M83").Value = -5764.25
$ws.Range("N83").Value = -27755.428

# Row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws.Range("H97").Value = 1700
$ws.Range("I97").Value = 1700
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1700
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1204
$ws.Range("N97").ClearContents()

# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 1365.4706
$ws.Range("I102").Value = 1263.3125
$ws.Range("K102").Value = 1263.3125
$ws.Range("M102").Value = 358.6875

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic | Boar Leather
$ws.Range("H46").Value = 710.3
$ws.Range("I46").Value = 710.3
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 710.3
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -522.3
$ws.Range("N46").ClearContents()

# Row 93: Hide to Go Seek | Gagana Leather
$ws.Range("H93").Value = 886.2105
$ws.Range("I93").Value = 836.3333
$ws.Range("J93").Value = 931.1
$ws.Range("K93").Value = 836.3333
$ws.Range("L93").Value = 931.1
$ws.Range("M93").Value = 411.6667
$ws.Range("N93").Value = -3427.1

# Row 104: Brace Yourselves | Gazelleskin Bracers of Fending
$ws.Range("H104").Value = 14511.667
$ws.Range("J104").Value = 14511.667
$ws.Range("L104").Value = 14511.667
$ws.Range("N104").Value = -21499.667

$ws = $wb.Worksheets.Item("WVR")
# Row 104: Brimming with Confidence | Twinsilk Turban of Aiming
$ws.Range("H104").Value = 11549.75
$ws.Range("J104").Value = 11549.75
$ws.Range("L104").Value = 11549.75
$ws.Range("N104").Value = -18537.75

# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 347078.56
$ws.Range("I122").Value = 590352.25
$ws.Range("J122").Value = 2440.8333
$ws.Range("K122").Value = 1771056.75
$ws.Range("L122").Value = 7322.499899999999
$ws.Range("M122").Value = -1768606.75
$ws.Range("N122").Value = -12222.4999
